$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-24 11:08:21"

# Sheet "zh-cn": Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-24 11:08:16"
$wsZhCn.Range("K2").Value = "2016-08-24 11:08:51"

# Sheet "de-de": Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-24 11:08:21"
$wsDeDe.Range("K2").Value = "2016-08-24 11:08:58"
